$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.363.32'
$ws.Range("E2").Value = '  +12.58%  '

$ws.Range("D3").Value = '1.825.59'
$ws.Range("E3").Value = '  +9.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.75'
$ws.Range("E5").Value = '  +4.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +4.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.57'
$ws.Range("E8").Value = '  +6.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.95'
$ws.Range("E9").Value = '  +6.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0673'
$ws.Range("E11").Value = '  +4.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0932'
$ws.Range("E12").Value = '  +2.93%  '

$ws.Range("D13").Value = '2.088.27'
$ws.Range("E13").Value = '  +9.13%  '

$ws.Range("D14").Value = '1.830.67'
$ws.Range("E14").Value = '  +9.53%  '

$ws.Range("E15").Value = '  +5.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.45'
$ws.Range("E16").Value = '  +2.65%  '

$ws.Range("D17").Value = '34.317.80'
$ws.Range("E17").Value = '  +12.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.30'
$ws.Range("E18").Value = '  +7.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.79'
$ws.Range("E19").Value = '  +4.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '258.39'
$ws.Range("E20").Value = '  +6.47%  '

$ws.Range("D21").Value = '0.0₃0752'
$ws.Range("E21").Value = '  +4.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("E23").Value = '  +6.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.36'
$ws.Range("E24").Value = '  +1.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E25").Value = '  +3.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.57'
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.64'
$ws.Range("E27").Value = '  +5.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.14'
$ws.Range("E28").Value = '  +6.99%  '

$ws.Range("E29").Value = '  +2.73%  '

$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.92'
$ws.Range("E31").Value = '  +13.07%  '

$ws.Range("E32").Value = '  +5.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.22'
$ws.Range("E33").Value = '  +6.60%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.58'
$ws.Range("E34").Value = '  +9.09%  '

$ws.Range("D35").Value = '1.552.81'
$ws.Range("E35").Value = '  +3.71%  '

$ws.Range("E36").Value = '  +1.56%  '

$ws.Range("B37").Value = 'MinaProtocolToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.36'
$ws.Range("E37").Value = '  +230.04%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.09'
$ws.Range("E38").Value = '  +6.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.642'
$ws.Range("E39").Value = '  +6.95%  '

$ws.Range("E40").Value = '  +6.77%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '84.75'
$ws.Range("E41").Value = '  +0.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.81'
$ws.Range("E42").Value = '  +5.32%  '

$ws.Range("E43").Value = '  +9.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.33'
$ws.Range("E44").Value = '  +1.56%  '

$ws.Range("E45").Value = '  +10.35%  '

$ws.Range("E46").Value = '  +5.46%  '

$ws.Range("E47").Value = '  +5.35%  '

$ws.Range("D48").Value = '1.990.27'
$ws.Range("E48").Value = '  +10.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.39'
$ws.Range("E49").Value = '  +27.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.80'
$ws.Range("E50").Value = '  +4.63%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  -0.04%  '
